$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H135").Value = 4718
$ws.Range("I135").Value = 4718
$ws.Range("J135").Value = 0
$ws.Range("K135").Value = 42462
$ws.Range("L135").Value = 0
$ws.Range("M135").Value = -39927
$ws.Range("N135").ClearContents()

$ws.Range("H138").Value = 1849.81
$ws.Range("I138").Value = 1239.0741
$ws.Range("J138").Value = 2075.6987
$ws.Range("K138").Value = 3717.2223
$ws.Range("L138").Value = 6227.0961
$ws.Range("M138").Value = 1422.7777
$ws.Range("N138").Value = -16507.0961

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3376.85
$ws.Range("I32").Value = 2603.173
$ws.Range("K32").Value = 2603.173
$ws.Range("M32").Value = -2316.173

$ws.Range("H42").Value = 17257.75
$ws.Range("J42").Value = 17257.75
$ws.Range("L42").Value = 17257.75
$ws.Range("N42").Value = -18229.75

$ws.Range("H74").Value = 9760.08
$ws.Range("I74").Value = 1787.8096
$ws.Range("K74").Value = 1787.8096
$ws.Range("M74").Value = -913.8096

$ws.Range("H77").Value = 9760.08
$ws.Range("I77").Value = 1787.8096
$ws.Range("K77").Value = 8939.048000000001
$ws.Range("M77").Value = -4571.048000000001

$ws.Range("H102").Value = 3013.3
$ws.Range("I102").Value = 3013.3
$ws.Range("K102").Value = 3013.3
$ws.Range("M102").Value = -1391.3

$ws.Range("H132").Value = 22659.416
$ws.Range("I132").Value = 10541.333
$ws.Range("K132").Value = 31623.999
$ws.Range("M132").Value = -29093.999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H64").Value = 4223
$ws.Range("J64").Value = 4760.1665
$ws.Range("L64").Value = 4760.1665
$ws.Range("N64").Value = -5210.1665

$ws.Range("H67").Value = 4223
$ws.Range("J67").Value = 4760.1665
$ws.Range("L67").Value = 4760.1665
$ws.Range("N67").Value = -6320.1665

$ws.Range("H86").Value = 2784.2
$ws.Range("I86").Value = 2941.75
$ws.Range("J86").Value = 2504.111
$ws.Range("K86").Value = 2941.75
$ws.Range("L86").Value = 2504.111
$ws.Range("M86").Value = -1818.75
$ws.Range("N86").Value = -4750.111

$ws.Range("H89").Value = 2784.2
$ws.Range("I89").Value = 2941.75
$ws.Range("J89").Value = 2504.111
$ws.Range("K89").Value = 14708.75
$ws.Range("L89").Value = 12520.555
$ws.Range("M89").Value = -9092.75
$ws.Range("N89").Value = -23752.555

$ws.Range("H107").Value = 1682.375
$ws.Range("I107").Value = 1537.0555
$ws.Range("K107").Value = 1537.0555
$ws.Range("M107").Value = 382.9445000000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H8").Value = 6100.8
$ws.Range("I8").Value = 5669.6665
$ws.Range("J8").Value = 6747.5
$ws.Range("K8").Value = 5669.6665
$ws.Range("L8").Value = 6747.5
$ws.Range("M8").Value = -5529.6665
$ws.Range("N8").Value = -7027.5

$ws.Range("H23").Value = 7399.5

$ws.Range("H25").Value = 6279.6665
$ws.Range("I25").Value = 6185.6
$ws.Range("J25").Value = 6750
$ws.Range("K25").Value = 6185.6
$ws.Range("L25").Value = 6750
$ws.Range("M25").Value = -6011.6
$ws.Range("N25").Value = -7098

$ws.Range("H27").Value = 7399.5

$ws.Range("H31").Value = 28604.316
$ws.Range("I31").Value = 19176
$ws.Range("K31").Value = 19176
$ws.Range("M31").Value = -18881

$ws.Range("H34").Value = 28604.316
$ws.Range("I34").Value = 19176
$ws.Range("K34").Value = 19176
$ws.Range("M34").Value = -18974

$ws.Range("H132").Value = 15940.75
$ws.Range("I132").Value = 3930.2856
$ws.Range("K132").Value = 11790.8568
$ws.Range("M132").Value = -9260.856800000001

$ws.Range("H137").Value = 64454.547
$ws.Range("J137").Value = 64454.547
$ws.Range("L137").Value = 64454.547
$ws.Range("N137").Value = -74654.54699999999

$ws.Range("H140").Value = 107864.445
$ws.Range("J140").Value = 107864.445
$ws.Range("L140").Value = 107864.445
$ws.Range("N140").Value = -118224.445

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 1002.9643
$ws.Range("I113").Value = 1038.7778
$ws.Range("J113").Value = 986
$ws.Range("K113").Value = 3116.3334
$ws.Range("L113").Value = 2958
$ws.Range("M113").Value = -946.3334000000004
$ws.Range("N113").Value = -7298

$ws.Range("H137").Value = 6498.3335
$ws.Range("I137").Value = 5747.5
$ws.Range("J137").Value = 8000
$ws.Range("K137").Value = 17242.5
$ws.Range("L137").Value = 24000
$ws.Range("M137").Value = -12142.5
$ws.Range("N137").Value = -34200

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H21").Value = 11179.25
$ws.Range("J21").Value = 12741.5
$ws.Range("L21").Value = 12741.5
$ws.Range("N21").Value = -13087.5

$ws.Range("H30").Value = 11179.25
$ws.Range("J30").Value = 12741.5
$ws.Range("L30").Value = 12741.5
$ws.Range("N30").Value = -12951.5

$ws.Range("H32").Value = 229998.33
$ws.Range("J32").Value = 229998.33
$ws.Range("L32").Value = 229998.33
$ws.Range("N32").Value = -230590.33

$ws.Range("H41").Value = 900
$ws.Range("I41").Value = 900
$ws.Range("K41").Value = 900
$ws.Range("M41").Value = -545

$ws.Range("H95").Value = 30344
$ws.Range("J95").Value = 30344
$ws.Range("L95").Value = 30344
$ws.Range("N95").Value = -35836

$ws.Range("H120").Value = 49090.816
$ws.Range("J120").Value = 49090.816
$ws.Range("L120").Value = 49090.816
$ws.Range("N120").Value = -58766.816

$ws.Range("H132").Value = 10696.25
$ws.Range("I132").Value = 8226.532999999999
$ws.Range("J132").Value = 23044.834
$ws.Range("K132").Value = 24679.599
$ws.Range("L132").Value = 69134.50199999999
$ws.Range("M132").Value = -22149.599
$ws.Range("N132").Value = -74194.50199999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 27029324
$ws.Range("I22").Value = 1875
$ws.Range("J22").Value = 62502852
$ws.Range("K22").Value = 1875
$ws.Range("L22").Value = 62502852
$ws.Range("M22").Value = -1580
$ws.Range("N22").Value = -62503442

$ws.Range("H23").Value = 10782.4
$ws.Range("I23").Value = 8137.3335
$ws.Range("J23").Value = 14750
$ws.Range("K23").Value = 8137.3335
$ws.Range("L23").Value = 14750
$ws.Range("M23").Value = -7907.3335
$ws.Range("N23").Value = -15210

$ws.Range("H25").Value = 983998.75
$ws.Range("I25").Value = 16000
$ws.Range("J25").Value = 1306665
$ws.Range("K25").Value = 16000
$ws.Range("L25").Value = 1306665
$ws.Range("M25").Value = -15770
$ws.Range("N25").Value = -1307125

$ws.Range("H27").Value = 27029324
$ws.Range("I27").Value = 1875
$ws.Range("J27").Value = 62502852
$ws.Range("K27").Value = 1875
$ws.Range("L27").Value = 62502852
$ws.Range("M27").Value = -1768
$ws.Range("N27").Value = -62503066

$ws.Range("H46").Value = 3848.4167
$ws.Range("I46").Value = 3597
$ws.Range("K46").Value = 3597
$ws.Range("M46").Value = -3409

$ws.Range("H122").Value = 23262336
$ws.Range("J122").Value = 6643.12
$ws.Range("L122").Value = 19929.36
$ws.Range("N122").Value = -24829.36

$ws.Range("H132").Value = 67802
$ws.Range("I132").Value = 3397
$ws.Range("K132").Value = 10191
$ws.Range("M132").Value = -7661

$ws.Range("H136").Value = 81562.39
$ws.Range("I136").Value = 202250.7
$ws.Range("J136").Value = 14513.333
$ws.Range("K136").Value = 606752.1000000001
$ws.Range("L136").Value = 43539.999
$ws.Range("M136").Value = -604202.1000000001
$ws.Range("N136").Value = -48639.999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H68").Value = 38500
$ws.Range("J68").Value = 38500
$ws.Range("L68").Value = 38500
$ws.Range("N68").Value = -40122

$ws.Range("H69").Value = 13144.429
$ws.Range("J69").Value = 13144.429
$ws.Range("L69").Value = 13144.429
$ws.Range("N69").Value = -14642.429

$ws.Range("H71").Value = 38500
$ws.Range("J71").Value = 38500
$ws.Range("L71").Value = 115500
$ws.Range("N71").Value = -123612

$ws.Range("H72").Value = 13144.429
$ws.Range("J72").Value = 13144.429
$ws.Range("L72").Value = 39433.287
$ws.Range("N72").Value = -46921.287

$ws.Range("H96").Value = 1633.8
$ws.Range("I96").Value = 1300
$ws.Range("J96").Value = 1755.1818
$ws.Range("K96").Value = 1300
$ws.Range("L96").Value = 1755.1818
$ws.Range("M96").Value = 73
$ws.Range("N96").Value = -4501.1818

$ws.Range("H113").Value = 2128.0667
$ws.Range("J113").Value = 955.875
$ws.Range("L113").Value = 2867.625
$ws.Range("N113").Value = -7207.625

$ws.Range("H122").Value = 521089.06
$ws.Range("I122").Value = 719999.7
$ws.Range("J122").Value = 7236.5835
$ws.Range("K122").Value = 2159999.1
$ws.Range("L122").Value = 21709.7505
$ws.Range("M122").Value = -2157549.1
$ws.Range("N122").Value = -26609.7505

$ws.Range("H132").Value = 6739.643
$ws.Range("I132").Value = 2140.7334
$ws.Range("K132").Value = 6422.2002
$ws.Range("M132").Value = -3892.2002

$ws.Range("H136").Value = 9793.59
$ws.Range("I136").Value = 1015.34485
$ws.Range("J136").Value = 35250.5
$ws.Range("K136").Value = 3046.03455
$ws.Range("L136").Value = 105751.5
$ws.Range("M136").Value = -496.0345499999999
$ws.Range("N136").Value = -110851.5
